$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").Value = 8.427199999999997
$ws.Range("A8").Value = -22.08450000000001
$ws.Range("A10").Value = -22.04380000000001
$ws.Range("A12").Value = -21.95810000000001
$ws.Range("B12").Value = 6.378199999999998
$ws.Range("C12").Value = -11.73290000000001
$ws.Range("C13").Value = -12.0113
$ws.Range("B15").Value = 5.585399999999995
$ws.Range("B17").Value = 5.0766
$ws.Range("A18").Value = -22.33030000000001
$ws.Range("C21").Value = -13.7626
$ws.Range("C25").Value = -12.3005
$ws.Range("B26").Value = 5.3892
$ws.Range("B27").Value = 6.0097
$ws.Range("B28").Value = 6.104199999999999
$ws.Range("C32").Value = -12.621
$ws.Range("C36").Value = -11.86110000000001
$ws.Range("A37").Value = -21.846
$ws.Range("B37").Value = 5.629
$ws.Range("C38").Value = -12.0309
$ws.Range("C41").Value = -12.94490000000001
$ws.Range("B47").Value = 6.0112
$ws.Range("C52").Value = -11.36570000000001
$ws.Range("A55").Value = -22.0577
$ws.Range("C59").Value = -12.77410000000001
$ws.Range("B65").Value = 6.253799999999997
$ws.Range("C67").Value = -11.83399999999999
$ws.Range("A68").Value = -21.4509
$ws.Range("B73").Value = 9.292999999999997
$ws.Range("A77").Value = -20.61539999999999
$ws.Range("A78").Value = -19.65719999999999
$ws.Range("A81").Value = -22.12230000000001
$ws.Range("A82").Value = -21.69580000000001
$ws.Range("B84").Value = 5.106899999999999
$ws.Range("C84").Value = -13.179
$ws.Range("B85").Value = 5.3247
$ws.Range("C88").Value = -13.0693
$ws.Range("C89").Value = -14.3341
$ws.Range("B93").Value = 5.471300000000002
$ws.Range("B95").Value = 6.376700000000005
$ws.Range("C95").Value = -12.7588
$ws.Range("B98").Value = 5.737000000000003
$ws.Range("B99").Value = 5.482699999999999
$ws.Range("B101").Value = 5.5658
$ws.Range("C105").Value = -12.72040000000001
